$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conditions")

# Manipulation check: change price strings from "cena: X zł" to "CENA: X zł"
$ws.Range("D2").Value = "CENA: 2259 zł"
$ws.Range("D3").Value = "CENA: 719 zł"
$ws.Range("D4").Value = "CENA: 1939 zł"
$ws.Range("D5").Value = "CENA: 249 zł"

# Update selection / view position
$ws.Range("D6").Select()
